function Add-ScrimRow {
    param($ws, $row, $tmplRow, $values)
    $src = $ws.Range("A" + $tmplRow + ":N" + $tmplRow)
    $src.Copy()
    $dst = $ws.Range("A" + $row + ":N" + $row)
    $dst.PasteSpecial(-4122)
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value2 = $values[$c]
    }
}

$wb = $excel.ActiveWorkbook

# ===== Sheet: Ring of Fire =====
$ws = $wb.Worksheets.Item("Ring of Fire")
$tmplS4 = 4
$tmplS5 = 8

Add-ScrimRow $ws 23 $tmplS4 @("JUJU", "CHARLIE", "PAM", "MEG", "BERRY", "JAE-YONG", "Equipo 2", "PLP|BrriN", "MTM|snoiy", "PLP|Mine", "NHG|Xemp", "NHG|Bayarea", "NHG|GN", "20250724T014508.000Z")
Add-ScrimRow $ws 24 $tmplS4 @("JUJU", "CHARLIE", "PAM", "MEG", "BERRY", "JAE-YONG", "Equipo 2", "PLP|BrriN", "MTM|snoiy", "PLP|Mine", "NHG|Xemp", "NHG|Bayarea", "NHG|GN", "20250724T014333.000Z")
Add-ScrimRow $ws 25 $tmplS4 @("DOUG", "BO", "AMBER", "LILY", "JAE-YONG", "BONNIE", "Equipo 2", "PLP|BrriN", "PLP|Mine", "MTM|snoiy", "NHG|Xemp", "NHG|GN", "NHG|Bayarea", "20250724T013841.000Z")
Add-ScrimRow $ws 26 $tmplS4 @("DOUG", "BO", "AMBER", "LILY", "JAE-YONG", "BONNIE", "Equipo 2", "PLP|BrriN", "PLP|Mine", "MTM|snoiy", "NHG|Xemp", "NHG|GN", "NHG|Bayarea", "20250724T013629.000Z")
Add-ScrimRow $ws 27 $tmplS4 @("BEA", "KENJI", "CROW", "ALLI", "LUMI", "JAE-YONG", "Equipo 2", "PLP|BrriN", "MTM|snoiy", "PLP|Mine", "NHG|Xemp", "NHG|Bayarea", "NHG|GN", "20250724T013157.000Z")
Add-ScrimRow $ws 28 $tmplS4 @("BEA", "KENJI", "CROW", "ALLI", "LUMI", "JAE-YONG", "Equipo 2", "PLP|BrriN", "MTM|snoiy", "PLP|Mine", "NHG|Xemp", "NHG|Bayarea", "NHG|GN", "20250724T012940.000Z")

# ===== Sheet: Pit Stop =====
$ws = $wb.Worksheets.Item("Pit Stop")
$tmplS4 = 6
$tmplS5 = 4

Add-ScrimRow $ws 28 $tmplS5 @("EL PRIMO", "MELODIE", "SHADE", "KAZE", "HANK", "LOU", "Equipo 1", "PLP|BrriN", "MTM|snoiy", "PLP|Mine", "NHG|Xemp", "NHG|GN", "NHG|Bayarea", "20250724T012247.000Z")

# ===== Sheet: Crystal Arcade =====
$ws = $wb.Worksheets.Item("Crystal Arcade")
$tmplS4 = 6
$tmplS5 = 4

Add-ScrimRow $ws 36 $tmplS5 @("DRACO", "BO", "JAE-YONG", "CROW", "MORTIS", "EMZ", "Equipo 1", "TE|Rafikii", "TE|Ezlivi", "TE|Belal", "TRB|Zeus 解開", "TRB|Lxffy", "TRB|R B M", "20250724T015301.000Z")
Add-ScrimRow $ws 37 $tmplS5 @("DRACO", "BO", "JAE-YONG", "CROW", "MORTIS", "EMZ", "Equipo 1", "TE|Rafikii", "TE|Ezlivi", "TE|Belal", "TRB|Zeus 解開", "TRB|Lxffy", "TRB|R B M", "20250724T015106.000Z")
Add-ScrimRow $ws 38 $tmplS4 @("DRACO", "BO", "JAE-YONG", "CROW", "MORTIS", "EMZ", "Equipo 2", "TE|Rafikii", "TE|Ezlivi", "TE|Belal", "TRB|Zeus 解開", "TRB|Lxffy", "TRB|R B M", "20250724T014927.000Z")
Add-ScrimRow $ws 39 $tmplS5 @("AMBER", "JANET", "KAZE", "SANDY", "LILY", "TARA", "Equipo 1", "TE|Rafikii", "TE|Ezlivi", "TE|Belal", "TRB|Zeus 解開", "TRB|Lxffy", "TRB|R B M", "20250724T014330.000Z")
Add-ScrimRow $ws 40 $tmplS5 @("AMBER", "JANET", "KAZE", "SANDY", "LILY", "TARA", "Equipo 1", "TE|Rafikii", "TE|Ezlivi", "TE|Belal", "TRB|Zeus 解開", "TRB|Lxffy", "TRB|R B M", "20250724T014150.000Z")

# ===== Sheet: New Horizons =====
$ws = $wb.Worksheets.Item("New Horizons")
$tmplS4 = 8
$tmplS5 = 4

Add-ScrimRow $ws 40 $tmplS5 @("CHARLIE", "KAZE", "SQUEAK", "GENE", "JAE-YONG", "CORDELIUS", "Equipo 1", "TE|Rafikii", "TE|Belal", "TE|Ezlivi", "TRB|Zeus 解開", "TRB|R B M", "TRB|Lxffy", "20250724T013435.000Z")
Add-ScrimRow $ws 41 $tmplS5 @("CHARLIE", "KAZE", "SQUEAK", "GENE", "JAE-YONG", "CORDELIUS", "Equipo 1", "TE|Rafikii", "TE|Belal", "TE|Ezlivi", "TRB|Zeus 解開", "TRB|R B M", "TRB|Lxffy", "20250724T013314.000Z")
Add-ScrimRow $ws 42 $tmplS4 @("BELLE", "GRAY", "ANGELO", "CORDELIUS", "CHUCK", "RUFFS", "Equipo 2", "TE|Rafikii", "TE|Ezlivi", "TE|Belal", "TRB|Lxffy", "TRB|R B M", "TRB|Zeus 解開", "20250724T012752.000Z")
Add-ScrimRow $ws 43 $tmplS4 @("BELLE", "GRAY", "ANGELO", "CORDELIUS", "CHUCK", "RUFFS", "Equipo 2", "TE|Rafikii", "TE|Ezlivi", "TE|Belal", "TRB|Lxffy", "TRB|R B M", "TRB|Zeus 解開", "20250724T012537.000Z")

$excel.CutCopyMode = $false
Write-Output "Done applying scrims update"
